# "Generate Report for Handoff" — refresh the localization-status report:
#   * zh-cn / de-de rows move from "Handed back: in sync with en-US" to
#     "Ready for handoff"
#   * the associated generate/handoff timestamps advance a few seconds
#   * the Status column (now holding a shorter label) is narrowed on each
#     per-locale sheet, matching the narrower "Overview" summary columns

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"               # zh-cn status
$ws.Range("F2").Value = "Ready for handoff"                # de-de status
$ws.Range("G2").Value = "2016-08-27 20:58:13"              # Latest HO Xliff Generate Date

$ws.Columns.Item(5).ColumnWidth = 16.33                    # zh-cn col
$ws.Columns.Item(6).ColumnWidth = 16.33                    # de-de col

# ---- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"                 # Status
$ws.Range("H2").Value = "2016-08-27 20:58:09"               # Latest Handoff Datetime

$ws.Columns.Item(3).ColumnWidth = 16.33                     # Status col

# ---- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"                 # Status
$ws.Range("H2").Value = "2016-08-27 20:58:13"                # Latest Handoff Datetime

$ws.Columns.Item(3).ColumnWidth = 16.33                      # Status col
